$d = $word.ActiveDocument

# Namespace declarations reused for each InsertXML payload.
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

# Map of exact (trimmed) paragraph text -> replacement w:p XML that reproduces
# the diff (bold "vs" headings, split runs + proofErr markers).
$replacements = @{
    "Accessible vs Inaccessible" = '<w:p ' + $wns + ' w14:paraId="780E95BC" w14:textId="390CBF54" w:rsidR="000C14CF" w:rsidRDefault="000C14CF"><w:pPr><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="000C14CF"><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>Accessible vs Inaccessible</w:t></w:r></w:p>'

    "Deterministic vs Non-Deterministic" = '<w:p ' + $wns + ' w14:paraId="2FB3AD91" w14:textId="39CC730F" w:rsidR="000C14CF" w:rsidRDefault="000C14CF"><w:pPr><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="000C14CF"><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Deterministic vs </w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidRPr="000C14CF"><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>Non-Deterministic</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>'

    "Episodic vs Non-Episodic" = '<w:p ' + $wns + ' w14:paraId="3596EAD8" w14:textId="540D8B98" w:rsidR="000C14CF" w:rsidRDefault="000C14CF"><w:pPr><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="000C14CF"><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Episodic vs </w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidRPr="000C14CF"><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>Non-Episodic</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>'

    "Static vs Dynamic" = '<w:p ' + $wns + ' w14:paraId="43903BDB" w14:textId="77DEDC44" w:rsidR="000C14CF" w:rsidRDefault="000C14CF"><w:pPr><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="000C14CF"><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>Static vs Dynamic</w:t></w:r></w:p>'

    "n(sbj, obj, ant, con)" = '<w:p ' + $wns + ' w14:paraId="7961EB29" w14:textId="58F4349A" w:rsidR="00D70518" w:rsidRDefault="00D70518"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r w:rsidRPr="00D70518"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>n(</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramEnd"/><w:r w:rsidRPr="00D70518"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>sbj</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00D70518"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>, obj, ant, con)</w:t></w:r></w:p>'

    "commitment, prohibition, authorization, power, sanction" = '<w:p ' + $wns + ' w14:paraId="56E5EA86" w14:textId="01712FA5" w:rsidR="00D70518" w:rsidRDefault="00D70518"><w:pPr><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="00D70518"><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>commitment, prohibition, authorization, power, sanctio</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>n</w:t></w:r></w:p>'

    "Informative Directive Commissive" = '<w:p ' + $wns + ' w14:paraId="2944EA9C" w14:textId="1FC84DFC" w:rsidR="00662EE5" w:rsidRDefault="00662EE5"><w:pPr><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="00662EE5"><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>Informative Directive Commissive</w:t></w:r></w:p>'
}

# Walk the paragraphs collection; match on the exact (trimmed) paragraph text
# so image/drawing paragraphs and blank paragraphs are left untouched.
foreach ($p in $d.Paragraphs) {
    $raw = $p.Range.Text
    $key = $raw.TrimEnd([char]13, [char]10, [char]7)
    if ($replacements.ContainsKey($key)) {
        $p.Range.InsertXML($replacements[$key]) | Out-Null
    }
}
